# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# These two sheets mirror the same event data, so the same row/value
# updates are applied to both.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 11315
    4  = 10618
    5  = 599
    7  = 758
    9  = 31
    12 = 10517
    13 = 3245
    18 = 95
    19 = 406
    21 = 10826
    23 = 18
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
